$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 19.07.2024"

# Row 6
$ws.Range("B6").Value = "21.07."
$ws.Range("C6").Value = "22.07."
$ws.Range("D6").Value = "PAYPAL VMRVBE"
$ws.Range("E6").Value = "22,67-"

# Row 7
$ws.Range("B7").Value = "23.07."
$ws.Range("C7").Value = "24.07."
$ws.Range("D7").Value = "KARTENZ./23.07 LIDL RO"
$ws.Range("E7").Value = "84,09-"

# Row 8
$ws.Range("B8").Value = "24.07."
$ws.Range("C8").Value = "25.07."
$ws.Range("D8").Value = "BURGER KING Mellrichstadt"
$ws.Range("E8").Value = "33,02-"

# Row 9
$ws.Range("B9").Value = "27.07."
$ws.Range("C9").Value = "28.07."
$ws.Range("D9").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E9").Value = "50,00-"

# Row 10 (previously empty)
$ws.Range("B10").Value = "31.07."
$ws.Range("C10").Value = "01.08."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 45481583"
$ws.Range("E10").Value = "39,99-"
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Row 11 (previously empty)
$ws.Range("B11").Value = "01.08."
$ws.Range("C11").Value = "02.08."
$ws.Range("D11").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E11").Value = "64,27-"
$ws.Range("E9").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 04.08.2024"
$ws.Range("E12").Value = "294,04-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 11.08.2024"
